$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "absorbance"
$ws.Range("F2").Formula = "=D2-E2"
$ws.Range("F3:F47").Formula = "=D3-E3"

$ws.Range("F2").Select()
